# Revised version: described search strategy + added Cvasciuc 2020 paper
#
# The study table on Sheet1 gains a new row for the "Cvasciuc 2020" paper
# (Northern Ireland, PH+PGL tumor), inserted in its sorted position
# (by year_pub) right after "Berends 2018" (row 10) and before
# "Ebbehoj 2020" (previously row 11, now row 12).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 11, pushing the existing rows 11-14 down
# to rows 12-15.
$ws.Rows.Item(11).Insert()

# Populate the new row with the new study's data.
$ws.Cells.Item(11, 2).Value() = "Cvasciuc 2020"     # B11 ref
$ws.Cells.Item(11, 3).Value() = 2020                # C11 year_pub
$ws.Cells.Item(11, 4).Value() = 2010                # D11 year_from
$ws.Cells.Item(11, 5).Value() = 2018                # E11 year_to
$ws.Cells.Item(11, 6).Value() = "Northern Ireland"  # F11 location
$ws.Cells.Item(11, 7).Value() = "PH+PGL"            # G11 tumor
$ws.Cells.Item(11, 8).Value() = 64                  # H11 altitude
$ws.Cells.Item(11, 9).Value() = 86                  # I11 total_cases

# The "Leung 2021" row (now row 15) also gets its "inletter" flag (column
# A) set to 1, matching the other updated rows.
$ws.Cells.Item(15, 1).Value() = 1

# Move the active selection, matching the author's final cursor position.
$ws.Range("A16").Select()

# Re-apply the AutoFilter so its range grows from A1:J18 to A1:J19 to
# cover the newly inserted row (toggle off, then back on with the new
# range).
$ws.Range("A1:J19").AutoFilter()
$ws.Range("A1:J19").AutoFilter()

# Keep the hidden _FilterDatabase defined name in sync with the new
# AutoFilter range.
$names = $wb.Names
for ($i = 1; $i -le $names.Count(); $i++) {
    $n = $names.Item($i)
    if ($n.Name() -like "*_FilterDatabase*") {
        $n.RefersTo() = "=Sheet1!`$A`$1:`$J`$19"
    }
}
